$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.526.34"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.484.98"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.95"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.12"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("E7").Value = "  -1.03%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  -1.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.58"
$ws.Range("E10").Value = "  -3.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0783"
$ws.Range("E11").Value = "  +0.09%  "

$ws.Range("E12").Value = "  +1.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.868.83"
$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("E14").Value = "  -2.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.67"
$ws.Range("E15").Value = "  +7.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.457.97"
$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.751"
$ws.Range("E17").Value = "  -4.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.546.07"
$ws.Range("E18").Value = "  +0.61%  "

$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0930"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.95"
$ws.Range("E21").Value = "  +4.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.24"
$ws.Range("E22").Value = "  -2.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.07"
$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.70"
$ws.Range("E24").Value = "  -3.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"

$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.91"
$ws.Range("E27").Value = "  +1.87%  "

$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.36"
$ws.Range("E30").Value = "  +0.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.42"
$ws.Range("E31").Value = "  +2.66%  "

$ws.Range("E32").Value = "  -2.82%  "

$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.15"
$ws.Range("E34").Value = "  +5.96%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").Value = "  -5.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.94"
$ws.Range("E37").Value = "  -2.68%  "

$ws.Range("E38").Value = "  +2.06%  "

$ws.Range("E39").Value = "  -3.01%  "

$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.12"
$ws.Range("E41").Value = "  -4.27%  "

$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.77"
$ws.Range("E43").Value = "  -6.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.962.24"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0284"
$ws.Range("E45").Value = "  -0.41%  "

$ws.Range("E46").Value = "  -3.31%  "

$ws.Range("E47").Value = "  +2.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.727.12"
$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.49"
$ws.Range("E49").Value = "  -0.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.52"
$ws.Range("E50").Value = "  -3.59%  "

$ws.Range("E51").Value = "  -3.42%  "
